# Paths.xlsx — "Merged Historic FD and Projected FD in Total FD; MultiIndex in 6."
#
# Row 22 (D only): the generic "Final Demand" folder path used for the merged
# FD source file is replaced by a path pointing at the new "Merged FD Projected"
# file.
# Row 23: now describes the historical-FD source ("Historical FD" / path to
# "Merged FD Historical").
# Row 24: what used to be row 23 ("History") shifts down, keeping the plain
# "Final Demand" folder path.
# Row 25 (new row): the total/merged FD output ("FD Total" / path to "Total FD").
#
# Shared-string append order below matches the author's edit order so the
# resulting sharedStrings.xml indices line up with the canonical diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D22").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand\Merged FD Projected"
$ws.Range("D25").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand\Total FD"
$ws.Range("A23").Value = "Historical FD"
$ws.Range("D23").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand\Merged FD Historical"
$ws.Range("A25").Value = "FD Total"

$ws.Range("A24").Value = "History"
$ws.Range("D24").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand"

# Update the frozen-pane view: scroll down a bit and move the selection to A26,
# matching the saved view state in the edited workbook.
$ws.Activate() | Out-Null
$ws.Range("A26").Select() | Out-Null

$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 20
    $win.ScrollColumn = 4
}
